$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set full target grid (A1:D7) reflecting the corrected file extensions / new row.
$ws.Range("A1").Value = "File_info"
$ws.Range("B1").Value = "File"
$ws.Range("C1").Value = "File"
$ws.Range("D1").Value = "File"

$ws.Range("A2").Value = "Filename"
$ws.Range("B2").Value = "mockdata_niche_2018.tsv"
$ws.Range("C2").Value = "mockdata_dighum_2018.csv"

$ws.Range("A3").Value = "Internal unique identifier"
$ws.Range("B3").Value = "EID"
$ws.Range("C3").Value = "UT"

$ws.Range("A4").Value = "Title"
$ws.Range("B4").Value = "Title"
$ws.Range("C4").Value = "TI"
$ws.Range("D4").Value = "Title"

$ws.Range("A5").Value = "ISSN"
$ws.Range("B5").Value = "ISSN"
$ws.Range("C5").Value = "SN"

$ws.Range("A6").Value = "DOI"
$ws.Range("B6").Value = "DOI"
$ws.Range("C6").Value = "DI"

$ws.Range("A7").Value = "Organization unit"
$ws.Range("B7").Value = "subject"
$ws.Range("C7").Value = "subject"

$ws.Range("B2").Select()
